$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / note updates -------------------------------------------------
# Order matters: new shared-strings are appended to the workbook's string
# table in the order they are first introduced, so we set these in the
# same sequence implied by the target file.
$ws.Range("G24").Value = "Jocelyn is in OOP class"
$ws.Range("G14").Value = "David is in Web class"
$ws.Range("D16").Value = "3:30pm CT"
$ws.Range("D20").Value = "4:00pm CT"
$ws.Range("D17").Value = "3:00pm CT"
$ws.Range("D18").Value = "4:30pm CT"

# --- Fill in missing Date/Time values for rows 18 and 20 -----------------
# Column C otherwise defaults to a plain left-aligned style; copy the date
# number format from an already-formatted date cell before writing values.
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 43811
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 43811

# --- Row 17: bring Bailey Cook's entry up to date with Date + Time -------
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 43811

$excel.CutCopyMode = $false

# --- Highlight the two rows that need attention (conflicting classes) ----
$ws.Range("A13:G13").Interior.Color = 65535
$ws.Range("A17:G17").Interior.Color = 65535

# --- Restore the saved selection -----------------------------------------
$ws.Range("C41").Select()
